$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a daily price table for "Melón" (Macroferia Regional de Talca).
# This edit inserts two new blocks of 4 daily records each into the middle of the
# table (pushing the existing rows below them down), growing the sheet from
# A1:R198 to A1:R206.

# --- Block 1: insert 4 new rows before current row 121 -----------------------
$ws.Rows("121:124").Insert()

$data1 = @(
    @(44554, "Tuna",     "Extra",   2000, 900,  900,  900,  900),
    @(44554, "Calameño", "Primera", 2000, 700,  700,  700,  700),
    @(44554, "Tuna",     "Extra",   2000, 900,  900,  900,  900),
    @(44554, "Tuna",     "Primera", 2000, 700,  700,  700,  700)
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = 121 + $i
    $row = $data1[$i]
    $ws.Cells.Item($r, 1).Value = 5
    $ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($r, 3).Value = "Maule"
    $ws.Cells.Item($r, 4).Value = $row[0]
    $ws.Cells.Item($r, 5).Value = 7
    $ws.Cells.Item($r, 6).Value = 100112027
    $ws.Cells.Item($r, 7).Value = "Melón"
    $ws.Cells.Item($r, 8).Value = $row[1]
    $ws.Cells.Item($r, 9).Value = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
    $ws.Cells.Item($r, 11).Value = $row[4]
    $ws.Cells.Item($r, 12).Value = $row[5]
    $ws.Cells.Item($r, 13).Value = $row[6]
    $ws.Cells.Item($r, 14).Value = '$/unidad'
    $ws.Cells.Item($r, 15).Value = "Región del Maule"
    $ws.Cells.Item($r, 16).Value = $row[7]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

# --- Block 2: insert 4 new rows before current row 147 -----------------------
$ws.Rows("147:150").Insert()

$data2 = @(
    @(44553, "Calameño", "Extra",   4000, 1000, 1000, 1000, 1000),
    @(44553, "Calameño", "Primera", 4000, 800,  800,  800,  800),
    @(44553, "Tuna",     "Extra",   3000, 1000, 1000, 1000, 1000),
    @(44553, "Tuna",     "Primera", 3000, 800,  800,  800,  800)
)

for ($i = 0; $i -lt $data2.Length; $i++) {
    $r = 147 + $i
    $row = $data2[$i]
    $ws.Cells.Item($r, 1).Value = 5
    $ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($r, 3).Value = "Maule"
    $ws.Cells.Item($r, 4).Value = $row[0]
    $ws.Cells.Item($r, 5).Value = 7
    $ws.Cells.Item($r, 6).Value = 100112027
    $ws.Cells.Item($r, 7).Value = "Melón"
    $ws.Cells.Item($r, 8).Value = $row[1]
    $ws.Cells.Item($r, 9).Value = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
    $ws.Cells.Item($r, 11).Value = $row[4]
    $ws.Cells.Item($r, 12).Value = $row[5]
    $ws.Cells.Item($r, 13).Value = $row[6]
    $ws.Cells.Item($r, 14).Value = '$/unidad'
    $ws.Cells.Item($r, 15).Value = "Región del Maule"
    $ws.Cells.Item($r, 16).Value = $row[7]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
